$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits (row numbers refer to the pre-sort layout) ---

# Carlos Penzini (Logistics / Member): Major "CS" -> "CS & Math"
$ws.Cells.Item(20, 5).Value = "CS & Math"

# Isabelle Xiong (Sponsorship / Member): Major "CS, Neuroscience & Econ" -> "CS"
$ws.Cells.Item(33, 5).Value = "CS"

# Amy Yin (Marketing / Co-Lead): add Major + LinkedIn hyperlink
$ws.Cells.Item(41, 5).Value = "Math & CS"
$ws.Hyperlinks.Add($ws.Cells.Item(41, 6), "https://www.linkedin.com/in/amy-heqing-yin-b20a73251/")

# Debi Ahitov (Marketing / Member): Major "Bio" -> "Bio & CS"
$ws.Cells.Item(42, 5).Value = "Bio & CS"

# Mariam Gvenetadze (was Director / Treasurer): Team "Director" -> "Treasury"
$ws.Cells.Item(44, 3).Value = "Treasury"

# --- Re-sort the member table by Team, Title, Last, First (ascending) ---
$dataRange = $ws.Range("A2:G44")
$dataRange.Sort($ws.Range("C2:C44"), 1, $ws.Range("D2:D44"), [System.Reflection.Missing]::Value, 1, $ws.Range("B2:B44"), 1, $ws.Range("A2:A44"))

# --- View state tweaks ---
$ws.Range("I25").Select()

Write-Output "done"
